$wb = $excel.ActiveWorkbook

# Worksheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 1172
$ws.Range("F5").Value = 72
$ws.Range("F6").Value = 482
$ws.Range("F7").Value = 834
$ws.Range("F8").Value = 436
$ws.Range("F10").Value = 2084
$ws.Range("F11").Value = 591
$ws.Range("F14").Value = 975
$ws.Range("F15").Value = 137
$ws.Range("F16").Value = 2109
$ws.Range("F17").Value = 581
$ws.Range("F18").Value = 10313
$ws.Range("F19").Value = 1014
$ws.Range("F21").Value = 100
$ws.Range("F22").Value = 124
$ws.Range("F24").Value = 251

# Worksheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("F10").Value = 138
$ws.Range("F11").Value = 9
$ws.Range("F13").Value = 41
$ws.Range("F15").Value = 2
$ws.Range("F25").Value = 28

# Worksheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 5653
$ws.Range("F3").Value = 455
$ws.Range("F4").Value = 435

# Worksheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 5653
$ws.Range("F4").Value = 455
$ws.Range("F5").Value = 435
$ws.Range("F7").Value = 1172
$ws.Range("F10").Value = 72
$ws.Range("F11").Value = 482
$ws.Range("F12").Value = 834
$ws.Range("F14").Value = 436
$ws.Range("F16").Value = 2084
$ws.Range("F17").Value = 591
$ws.Range("F22").Value = 975
$ws.Range("F24").Value = 137
$ws.Range("F25").Value = 138
$ws.Range("F26").Value = 9
$ws.Range("F27").Value = 2109
$ws.Range("F28").Value = 581
$ws.Range("F30").Value = 41
$ws.Range("F31").Value = 1014
$ws.Range("F33").Value = 100
$ws.Range("F34").Value = 124
$ws.Range("F37").Value = 2
$ws.Range("F39").Value = 251
$ws.Range("F49").Value = 28
